# Generate Report for Archive
#
# The "Status" columns on every sheet (Overview!E:F = zh-cn/de-de status,
# and the Status column on the per-locale sheets) flip from
# "Ready for handoff" to "In Translation", and those now-narrower columns
# are resized to fit the shorter text.

$wb = $excel.ActiveWorkbook

# 1) Replace the status text everywhere it appears, on every worksheet.
#    NOTE: the literal must be the left operand of -eq — PowerShell's -eq
#    coerces the right operand to the left operand's type, and a boolean
#    cell's $true would otherwise "match" any non-empty string literal.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Value()) {
            $cell.Value = "In Translation"
        }
    }
}

# 2) Re-fit the columns that held that text so they match the shorter value.
#    Overview: Status-per-locale columns are E (zh-cn) and F (de-de).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

#    zh-cn / de-de: Status column is C.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
